$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "展览" (sheet1) - bump the "want-to-go" visitor counts (column F)
# ---------------------------------------------------------------------------
$wsA = $wb.Worksheets.Item("展览")
$wsA.Range("F2").Value  = 8184
$wsA.Range("F3").Value  = 127
$wsA.Range("F4").Value  = 97
$wsA.Range("F5").Value  = 32559
$wsA.Range("F12").Value = 811
$wsA.Range("F13").Value = 60
$wsA.Range("F14").Value = 620
$wsA.Range("F15").Value = 415
$wsA.Range("F17").Value = 564
$wsA.Range("F18").Value = 155
$wsA.Range("F19").Value = 424
$wsA.Range("F20").Value = 424
$wsA.Range("F23").Value = 727
$wsA.Range("F24").Value = 2375
$wsA.Range("F25").Value = 850
$wsA.Range("F26").Value = 75
$wsA.Range("F27").Value = 1099
$wsA.Range("F29").Value = 644
$wsA.Range("F30").Value = 7
$wsA.Range("F31").Value = 1086

# ---------------------------------------------------------------------------
# Sheet "演出" (sheet2) - same kind of counter bump
# ---------------------------------------------------------------------------
$wsB = $wb.Worksheets.Item("演出")
$wsB.Range("F3").Value  = 69
$wsB.Range("F4").Value  = 350
$wsB.Range("F10").Value = 1

# ---------------------------------------------------------------------------
# Sheet "全部类型" (sheet4) - counter bumps for the rows that are not touched
# structurally
# ---------------------------------------------------------------------------
$wsD = $wb.Worksheets.Item("全部类型")
$wsD.Range("F3").Value  = 8184
$wsD.Range("F4").Value  = 127
$wsD.Range("F5").Value  = 97
$wsD.Range("F7").Value  = 32560
$wsD.Range("F12").Value = 69
$wsD.Range("F15").Value = 350
$wsD.Range("F18").Value = 811
$wsD.Range("F19").Value = 60
$wsD.Range("F20").Value = 620
$wsD.Range("F21").Value = 415
$wsD.Range("F26").Value = 1
$wsD.Range("F28").Value = 564
$wsD.Range("F29").Value = 155
$wsD.Range("F30").Value = 424
$wsD.Range("F31").Value = 424
$wsD.Range("F34").Value = 727
$wsD.Range("F35").Value = 2375
$wsD.Range("F36").Value = 850
$wsD.Range("F37").Value = 75
$wsD.Range("F38").Value = 1099

# Row 41 ("Look Look") itself is untouched, but its event info now also
# occupies row 42 (with an updated counter) instead of the old "AP" entry
# that used to live there. Overwrite row 42's content in place.
$wsD.Range("F41").Value = 644

$wsD.Range("C42").Value = "广州·Look Look动漫嘉年华"
$wsD.Range("D42").Value = "东沙大道16号 健康方舟"
$wsD.Range("E42").Value = "2024.06.01 10:00-06.02 17:30"
$wsD.Range("F42").Value = 644
$wsD.Range("G42").Value = 29.9
$wsD.Range("H42").Value = "https://show.bilibili.com/platform/detail.html?id=82319"
$wsD.Range("I42").Value = "//i2.hdslb.com/bfs/openplatform/202403/Zv7tuBjf1709620427087.png"

# Insert a brand new row at position 43 for the "AP" event (pushes the old
# row 43 - the 622 volleyball event - down to row 44).
$wsD.Rows.Item(43).Insert()

# Give the new A43 the same bordered/centered "index" style used by the
# other rows in column A, then set its index value.
$wsD.Range("A42").Copy()
$wsD.Range("A43").PasteSpecial(-4122)
$wsD.Range("A43").Value = 42

# B43 holds a literal "YYYY-MM-DD" string (not a real date) in this sheet,
# so force text formatting before assigning it, then drop the format again
# so no stray numeric format sticks to the cell.
$wsD.Range("B43").NumberFormat = "@"
$wsD.Range("B43").Value = "2024-06-01"
$wsD.Range("B43").ClearFormats()

$wsD.Range("C43").Value = "广州·第五届AP动漫嘉年华"
$wsD.Range("D43").Value = "西环路1号 广州岭南会展中心"
$wsD.Range("E43").Value = "2024.06.01 10:00-06.01 17:00"
$wsD.Range("F43").Value = 7
$wsD.Range("G43").Value = 55
$wsD.Range("H43").Value = "https://show.bilibili.com/platform/detail.html?id=83462"
$wsD.Range("I43").Value = "//i1.hdslb.com/bfs/openplatform/202403/ZR2jKMOg1711076939687.jpeg"

# Row 44 now holds what used to be row 43 (622 volleyball). Its serial
# index needs to move from 42 to 43, and its counter bumps too.
$wsD.Range("A44").Value = 43
$wsD.Range("F44").Value = 1086
